$d = $word.ActiveDocument

# Change 1: "6 frascos" -> "7 frascos" (in context ". Por exemplo, caso a fila tenha 6 frascos")
$d.Content.Find.Execute("caso a fila tenha 6 frascos", $true, $false, $false, $false, $false,
                         $true, 1, $false, "caso a fila tenha 7 frascos", 2)

# Change 2: merge runs "A entrada contém vários casos de teste" + ", cada caso é expresso "
# into a single run "A entrada contém vários casos de teste, cada caso é expresso "
$d.Content.Find.Execute("A entrada contém vários casos de teste, cada caso é expresso ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "A entrada contém vários casos de teste, cada caso é expresso ", 2)
